# Deploying to gh-pages from @ IDGRLP/Tumorkonferenzen-IG -- apply the
# metadata refresh that shipped with that build:
#   1. Rename the "Include from Multiple PrimärT" sheet to "Include #0".
#   2. Insert a new "Jurisdiction" metadata row (empty value) right after
#      "Contact" and before "Description" in the Metadata sheet - this
#      pushes Description/Purpose/Copyright/Immutable down by one row.
#   3. Refresh the "Date" metadata value to the new build timestamp.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item(2)

# --- 1. Rename the second sheet -------------------------------------------
$wsInclude.Name = "Include #0"

# --- 2. Insert the new "Jurisdiction" row ----------------------------------
# Insert a blank row at 11 (pushing "Description" and everything below it
# down by one) and give it the same look as the other property rows by
# copying the formatting from the row above ("Contact").
$wsMeta.Rows.Item(11).Insert()
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsMeta.Range("A11").Value = "Jurisdiction"

# Write an explicit (empty) text value into B11 - same text type as the
# sheet's other empty-value property cells - then restore the plain
# property-row formatting (the quote-entry above nudges the cell into a
# text style we don't want to keep).
$wsMeta.Range("B11").Value = "'"
$wsMeta.Range("B10").Copy()
$wsMeta.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Refresh the Date value ---------------------------------------------
$wsMeta.Range("B8").Value = "2024-09-17T19:55:11+00:00"
